# Weekly update: a new price observation is inserted as row 101 (pushing the
# existing rows 101-127 down to 102-128) for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Albahaca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 101, shifting rows 101:127
# down to 102:128 (keeping their data/format intact).
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly observation.
$ws.Range("A101").Value() = 4
$ws.Range("B101").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C101").Value() = "Los Lagos"
$ws.Range("D101").Value() = 44736
$ws.Range("E101").Value() = 10
$ws.Range("F101").Value() = 100112052
$ws.Range("G101").Value() = "Albahaca"
$ws.Range("H101").Value() = "Sin especificar"
$ws.Range("I101").Value() = "Primera"
$ws.Range("J101").Value() = 150
$ws.Range("K101").Value() = 6000
$ws.Range("L101").Value() = 7000
$ws.Range("M101").Value() = 6533
$ws.Range("N101").Value() = "`$/paquete"
$ws.Range("O101").Value() = "Región de Arica y Parinacota"
$ws.Range("P101").Value() = 6533
$ws.Range("Q101").Value() = 1
$ws.Range("R101").Value() = "Hortaliza"
